# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the refreshed data snapshot (per commit "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets index 1 / sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 634
$wsExhibit.Range("F3").Value = 209
$wsExhibit.Range("F4").Value = 668
$wsExhibit.Range("F5").Value = 573
$wsExhibit.Range("F6").Value = 318
$wsExhibit.Range("F7").Value = 2818
$wsExhibit.Range("F9").Value = 7972
$wsExhibit.Range("F10").Value = 205
$wsExhibit.Range("F13").Value = 388

# Sheet "全部类型" (Worksheets index 4 / sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 634
$wsAll.Range("F3").Value = 209
$wsAll.Range("F4").Value = 668
$wsAll.Range("F5").Value = 573
$wsAll.Range("F6").Value = 318
$wsAll.Range("F9").Value = 2818
$wsAll.Range("F11").Value = 7972
$wsAll.Range("F12").Value = 205
$wsAll.Range("F17").Value = 388
